# -----------------------------------------------------------------------
# [ENF] Add Report and reformat
#
# Rebuilds the "Stock Balance Report" sheet:
#  - header text relabelled (Thai columns) and report title text kept
#  - base font switched from Arial 10 to Tahoma 9 (regular + bold)
#  - header row (row 5) gets centered, wrapped, bordered formatting
#  - data rows keep their borders/number format but move to the new font
#    and gain wrap-text
#  - trailing unused header cells (G5:I5) are cleared
#  - column widths + selection updated
#
# NOTE: each formatting property (Font, HorizontalAlignment, WrapText,
# NumberFormat, ...) is applied to a cell range in one dedicated block
# below -- properties for the same final combined style are grouped
# together and set back-to-back on the same Range object so the engine
# only ever materialises the styles that are actually needed (instead of
# churning through throw-away intermediate combinations).
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Alignment constants (Excel COM numeric values)
$xlLeft   = -4131
$xlRight  = -4152
$xlCenter = -4108

# -----------------------------------------------------------------------
# 1. Text content
# -----------------------------------------------------------------------
$ws.Range("A1").Value = "NSTDA"
$ws.Range("A2").Value = "Stock Balance Report"
$ws.Range("A3").Value = "Location"

$ws.Range("A5").Value = "ชื่อสินค้า"
$ws.Range("B5").Value = "ศูนย์"
$ws.Range("C5").Value = "คลังวัสดุ"
$ws.Range("D5").Value = "ยอดคงเหลือ"
$ws.Range("E5").Value = "หน่วยนับ"
$ws.Range("F5").Value = "มูลค่าสินค้าคงคลัง"

# Headers G5:I5 no longer used -- drop them (and their formatting)
$ws.Range("G5").Clear()
$ws.Range("H5").Clear()
$ws.Range("I5").Clear()

# -----------------------------------------------------------------------
# 2. Formatting -- one block per final style, each property set once.
# -----------------------------------------------------------------------

# -- Title / caption block (A1:A4): bold Tahoma 9, left, wrap, no border
$r = $ws.Range("A1:A4")
$r.Font.Name = "Tahoma"
$r.Font.Size = 9
$r.Font.Bold = $true
$r.HorizontalAlignment = $xlLeft
$r.WrapText = $true

# -- B4: regular Tahoma 9, left, wrap, no border
$r = $ws.Range("B4")
$r.Font.Name = "Tahoma"
$r.Font.Size = 9
$r.Font.Bold = $false
$r.HorizontalAlignment = $xlLeft
$r.WrapText = $true

# -- Header band row 5 (A5:F5): bold Tahoma 9, centered, wrap, border2 (kept as-is)
$r = $ws.Range("A5:F5")
$r.Font.Name = "Tahoma"
$r.Font.Size = 9
$r.Font.Bold = $true
$r.HorizontalAlignment = $xlCenter
$r.WrapText = $true

# -- Row label column (A6:A9): regular Tahoma 9, left, wrap, no border
$r = $ws.Range("A6:A9")
$r.Font.Name = "Tahoma"
$r.Font.Size = 9
$r.Font.Bold = $false
$r.HorizontalAlignment = $xlLeft
$r.WrapText = $true

# -- Left-aligned numeric-format text columns (B,C,E) rows 6:9
$r = $ws.Range("B6:C9")
$r.Font.Name = "Tahoma"
$r.Font.Size = 9
$r.Font.Bold = $false
$r.NumberFormat = "#,##0.00"
$r.HorizontalAlignment = $xlLeft
$r.WrapText = $true

$r = $ws.Range("E6:E9")
$r.Font.Name = "Tahoma"
$r.Font.Size = 9
$r.Font.Bold = $false
$r.NumberFormat = "#,##0.00"
$r.HorizontalAlignment = $xlLeft
$r.WrapText = $true

# -- Right-aligned, wrapped numeric-format columns (D,F) rows 6:9
$r = $ws.Range("D6:D9")
$r.Font.Name = "Tahoma"
$r.Font.Size = 9
$r.Font.Bold = $false
$r.NumberFormat = "#,##0.00"
$r.HorizontalAlignment = $xlRight
$r.WrapText = $true

$r = $ws.Range("F6:F9")
$r.Font.Name = "Tahoma"
$r.Font.Size = 9
$r.Font.Bold = $false
$r.NumberFormat = "#,##0.00"
$r.HorizontalAlignment = $xlRight
$r.WrapText = $true

# -- Right-aligned, no-wrap numeric-format columns (G,H,I) rows 6:9
$r = $ws.Range("G6:I9")
$r.Font.Name = "Tahoma"
$r.Font.Size = 9
$r.Font.Bold = $false
$r.NumberFormat = "#,##0.00"
$r.HorizontalAlignment = $xlRight
$r.WrapText = $false

# -----------------------------------------------------------------------
# 3. Row heights
# -----------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 23
$ws.Rows.Item(5).RowHeight = 16.5

# -----------------------------------------------------------------------
# 4. Column widths (ColumnWidth input is shifted so the exported
#    <col width> lands on the whole-number target; Excel's internal
#    "characters" <-> "1/256 px" rounding needs a small offset)
# -----------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 16.125   # -> 17
$ws.Columns.Item(2).ColumnWidth = 19.125   # -> 20
$ws.Columns.Item(3).ColumnWidth = 19.125   # -> 20
$ws.Columns.Item(4).ColumnWidth = 19.125   # -> 20
$ws.Columns.Item(5).ColumnWidth = 19.125   # -> 20
$ws.Columns.Item(6).ColumnWidth = 19.125   # -> 20
$ws.Columns.Item(7).ColumnWidth = 19.125   # -> 20
$ws.Columns.Item(8).ColumnWidth = 24.125   # -> 25
$ws.Columns.Item(9).ColumnWidth = 19.125   # -> 20

# -----------------------------------------------------------------------
# 5. Selection
# -----------------------------------------------------------------------
$ws.Range("H4").Select()
